$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 107.5
$ws.Range("I12").Value = 125.8
$ws.Range("J12").Value = 16
$ws.Range("K12").Value = 125.8
$ws.Range("L12").Value = 16
$ws.Range("M12").Value = 44.2
$ws.Range("N12").Value = -356
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("H32").Value = 2999.5
$ws.Range("I32").Value = 2999.5
$ws.Range("K32").Value = 2999.5
$ws.Range("M32").Value = -2673.5
$ws.Range("H33").Value = 1571.3846
$ws.Range("I33").Value = 558.8889
$ws.Range("J33").Value = 3849.5
$ws.Range("K33").Value = 558.8889
$ws.Range("L33").Value = 3849.5
$ws.Range("M33").Value = -329.8889
$ws.Range("N33").Value = -4307.5
$ws.Range("H40").Value = 27277486
$ws.Range("J40").Value = 37503492
$ws.Range("L40").Value = 37503492
$ws.Range("N40").Value = -37503842
$ws.Range("H51").Value = 9892.296
$ws.Range("I51").Value = 16331.333
$ws.Range("K51").Value = 16331.333
$ws.Range("M51").Value = -15847.333
$ws.Range("H55").Value = 306
$ws.Range("I55").Value = 87.09999999999999
$ws.Range("J55").Value = 743.8
$ws.Range("K55").Value = 87.09999999999999
$ws.Range("L55").Value = 743.8
$ws.Range("M55").Value = 126.9
$ws.Range("N55").Value = -1171.8
$ws.Range("H62").Value = 22734780
$ws.Range("I62").Value = 41672932
$ws.Range("K62").Value = 41672932
$ws.Range("M62").Value = -41672308
$ws.Range("H65").Value = 22734780
$ws.Range("I65").Value = 41672932
$ws.Range("K65").Value = 208364660
$ws.Range("M65").Value = -208361540
$ws.Range("H94").Value = 899.5
$ws.Range("I94").Value = 899.5
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 899.5
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -448.5
$ws.Range("N94").Value = $null
$ws.Range("H107").Value = 35717630
$ws.Range("I107").Value = 27781436
$ws.Range("K107").Value = 27781436
$ws.Range("M107").Value = -27779516
$ws.Range("H132").Value = 4987.2666
$ws.Range("I132").Value = 2842.037
$ws.Range("K132").Value = 8526.110999999999
$ws.Range("M132").Value = -5996.110999999999
$ws.Range("H133").Value = 99999.336
$ws.Range("J133").Value = 99999.336
$ws.Range("L133").Value = 99999.336
$ws.Range("N133").Value = -110119.336
$ws.Range("H136").Value = 434995
$ws.Range("J136").Value = 434995
$ws.Range("L136").Value = 434995
$ws.Range("N136").Value = -445195
$ws.Range("H137").Value = 44622.48
$ws.Range("I137").Value = 66724.35000000001
$ws.Range("J137").Value = 7049.3
$ws.Range("K137").Value = 200173.05
$ws.Range("L137").Value = 21147.9
$ws.Range("M137").Value = -197623.05
$ws.Range("N137").Value = -26247.9
$ws.Range("H138").Value = 3361.7
$ws.Range("I138").Value = 1602.4
$ws.Range("J138").Value = 3948.1333
$ws.Range("K138").Value = 4807.200000000001
$ws.Range("L138").Value = 11844.3999
$ws.Range("M138").Value = 332.7999999999993
$ws.Range("N138").Value = -22124.3999
$ws.Range("H141").Value = 3454.3333
$ws.Range("I141").Value = 3289.5625
$ws.Range("K141").Value = 9868.6875
$ws.Range("M141").Value = -4688.6875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1725.45
$ws.Range("I2").Value = 1333.2222
$ws.Range("J2").Value = 5255.5
$ws.Range("K2").Value = 1333.2222
$ws.Range("L2").Value = 5255.5
$ws.Range("M2").Value = -1220.2222
$ws.Range("N2").Value = -5481.5
$ws.Range("H32").Value = 9827.708000000001
$ws.Range("I32").Value = 6540.3486
$ws.Range("J32").Value = 38099
$ws.Range("K32").Value = 6540.3486
$ws.Range("L32").Value = 38099
$ws.Range("M32").Value = -6253.3486
$ws.Range("N32").Value = -38673
$ws.Range("H45").Value = 2520.52
$ws.Range("I45").Value = 2492.5
$ws.Range("K45").Value = 2492.5
$ws.Range("M45").Value = -2115.5
$ws.Range("H61").Value = 5166.25
$ws.Range("I61").Value = 1952.3846
$ws.Range("J61").Value = 13522.3
$ws.Range("K61").Value = 1952.3846
$ws.Range("L61").Value = 13522.3
$ws.Range("M61").Value = -1740.3846
$ws.Range("N61").Value = -13946.3
$ws.Range("H63").Value = 4022.3572
$ws.Range("I63").Value = 3622
$ws.Range("J63").Value = 4743
$ws.Range("K63").Value = 3622
$ws.Range("L63").Value = 4743
$ws.Range("M63").Value = -2936
$ws.Range("N63").Value = -6115
$ws.Range("H66").Value = 4022.3572
$ws.Range("I66").Value = 3622
$ws.Range("J66").Value = 4743
$ws.Range("K66").Value = 18110
$ws.Range("L66").Value = 23715
$ws.Range("M66").Value = -14678
$ws.Range("N66").Value = -30579
$ws.Range("H74").Value = 62871.125
$ws.Range("I74").Value = 69006.8
$ws.Range("J74").Value = 3559.6667
$ws.Range("K74").Value = 69006.8
$ws.Range("L74").Value = 3559.6667
$ws.Range("M74").Value = -68132.8
$ws.Range("N74").Value = -5307.6667
$ws.Range("H77").Value = 62871.125
$ws.Range("I77").Value = 69006.8
$ws.Range("J77").Value = 3559.6667
$ws.Range("K77").Value = 345034
$ws.Range("L77").Value = 17798.3335
$ws.Range("M77").Value = -340666
$ws.Range("N77").Value = -26534.3335
$ws.Range("H110").Value = 9285.15
$ws.Range("I110").Value = 9669.23
$ws.Range("J110").Value = 8571.857
$ws.Range("K110").Value = 9669.23
$ws.Range("L110").Value = 8571.857
$ws.Range("M110").Value = -7624.23
$ws.Range("N110").Value = -12661.857
$ws.Range("H116").Value = 1725.45
$ws.Range("I116").Value = 1333.2222
$ws.Range("J116").Value = 5255.5
$ws.Range("K116").Value = 1333.2222
$ws.Range("L116").Value = 5255.5
$ws.Range("M116").Value = 960.7778000000001
$ws.Range("N116").Value = -9843.5
$ws.Range("H122").Value = 2349.3076
$ws.Range("I122").Value = 2461.652
$ws.Range("J122").Value = 1488
$ws.Range("K122").Value = 7384.956
$ws.Range("L122").Value = 4464
$ws.Range("M122").Value = -4934.956
$ws.Range("N122").Value = -9364
$ws.Range("H132").Value = 2808.5454
$ws.Range("I132").Value = 3900.4736
$ws.Range("J132").Value = 1978.68
$ws.Range("K132").Value = 11701.4208
$ws.Range("L132").Value = 5936.04
$ws.Range("M132").Value = -9171.4208
$ws.Range("N132").Value = -10996.04
$ws.Range("H136").Value = 5166.25
$ws.Range("I136").Value = 1952.3846
$ws.Range("J136").Value = 13522.3
$ws.Range("K136").Value = 5857.1538
$ws.Range("L136").Value = 40566.89999999999
$ws.Range("M136").Value = -3307.1538
$ws.Range("N136").Value = -45666.89999999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1725.45
$ws.Range("I3").Value = 1333.2222
$ws.Range("J3").Value = 5255.5
$ws.Range("K3").Value = 1333.2222
$ws.Range("L3").Value = 5255.5
$ws.Range("M3").Value = -1219.2222
$ws.Range("N3").Value = -5483.5
$ws.Range("H20").Value = 2948.2
$ws.Range("I20").Value = 1271.4445
$ws.Range("J20").Value = 4320.091
$ws.Range("K20").Value = 1271.4445
$ws.Range("L20").Value = 4320.091
$ws.Range("M20").Value = -1024.4445
$ws.Range("N20").Value = -4814.091
$ws.Range("H22").Value = 417
$ws.Range("I22").Value = 417
$ws.Range("K22").Value = 417
$ws.Range("M22").Value = -244
$ws.Range("H99").Value = 4978
$ws.Range("I99").Value = 5281.5
$ws.Range("K99").Value = 5281.5
$ws.Range("M99").Value = -3783.5
$ws.Range("H105").Value = 1935.9131
$ws.Range("I105").Value = 1767.3684
$ws.Range("J105").Value = 2736.5
$ws.Range("K105").Value = 1767.3684
$ws.Range("L105").Value = 2736.5
$ws.Range("M105").Value = -20.36840000000007
$ws.Range("N105").Value = -6230.5
$ws.Range("H107").Value = 3659.5
$ws.Range("J107").Value = 3319.6
$ws.Range("L107").Value = 3319.6
$ws.Range("N107").Value = -7159.6
$ws.Range("H134").Value = 2381.3171
$ws.Range("I134").Value = 1800.7142
$ws.Range("K134").Value = 5402.142599999999
$ws.Range("M134").Value = -2867.142599999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4532.8184
$ws.Range("I16").Value = 1210.6666
$ws.Range("J16").Value = 8519.4
$ws.Range("K16").Value = 1210.6666
$ws.Range("L16").Value = 8519.4
$ws.Range("M16").Value = -923.6666
$ws.Range("N16").Value = -9093.4
$ws.Range("H31").Value = 214665.42
$ws.Range("I31").Value = 371495.62
$ws.Range("J31").Value = 2944.65
$ws.Range("K31").Value = 371495.62
$ws.Range("L31").Value = 2944.65
$ws.Range("M31").Value = -371200.62
$ws.Range("N31").Value = -3534.65
$ws.Range("H34").Value = 214665.42
$ws.Range("I34").Value = 371495.62
$ws.Range("J34").Value = 2944.65
$ws.Range("K34").Value = 371495.62
$ws.Range("L34").Value = 2944.65
$ws.Range("M34").Value = -371293.62
$ws.Range("N34").Value = -3348.65
$ws.Range("H58").Value = 3507.2285
$ws.Range("I58").Value = 3129.5
$ws.Range("K58").Value = 3129.5
$ws.Range("M58").Value = -2926.5
$ws.Range("H94").Value = 888.3043
$ws.Range("I94").Value = 614.9091
$ws.Range("K94").Value = 614.9091
$ws.Range("M94").Value = -163.9091
$ws.Range("H99").Value = 428297.53
$ws.Range("J99").Value = 15847.077
$ws.Range("L99").Value = 15847.077
$ws.Range("N99").Value = -18843.077
$ws.Range("H105").Value = 4269.413
$ws.Range("I105").Value = 1280.0435
$ws.Range("K105").Value = 1280.0435
$ws.Range("M105").Value = 466.9565
$ws.Range("H107").Value = 4701.108
$ws.Range("I107").Value = 761.2727
$ws.Range("K107").Value = 761.2727
$ws.Range("M107").Value = 1158.7273
$ws.Range("H113").Value = 4532.8184
$ws.Range("I113").Value = 1210.6666
$ws.Range("J113").Value = 8519.4
$ws.Range("K113").Value = 1210.6666
$ws.Range("L113").Value = 8519.4
$ws.Range("M113").Value = 959.3334
$ws.Range("N113").Value = -12859.4
$ws.Range("H122").Value = 2446.653
$ws.Range("I122").Value = 2349.3225
$ws.Range("J122").Value = 2614.2778
$ws.Range("K122").Value = 7047.967500000001
$ws.Range("L122").Value = 7842.8334
$ws.Range("M122").Value = -4597.967500000001
$ws.Range("N122").Value = -12742.8334
$ws.Range("H126").Value = 428297.53
$ws.Range("J126").Value = 15847.077
$ws.Range("L126").Value = 47541.231
$ws.Range("N126").Value = -52481.231
$ws.Range("H132").Value = 4376
$ws.Range("I132").Value = 1875.8889
$ws.Range("J132").Value = 6421.5454
$ws.Range("K132").Value = 5627.6667
$ws.Range("L132").Value = 19264.6362
$ws.Range("M132").Value = -3097.6667
$ws.Range("N132").Value = -24324.6362
$ws.Range("H134").Value = 2788.5757
$ws.Range("I134").Value = 2580.0356
$ws.Range("K134").Value = 7740.1068
$ws.Range("M134").Value = -5205.1068
$ws.Range("H136").Value = 3507.2285
$ws.Range("I136").Value = 3129.5
$ws.Range("K136").Value = 9388.5
$ws.Range("M136").Value = -6838.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 65.76470999999999
$ws.Range("I2").Value = 40.88889
$ws.Range("J2").Value = 93.75
$ws.Range("K2").Value = 245.33334
$ws.Range("L2").Value = 562.5
$ws.Range("M2").Value = -132.33334
$ws.Range("N2").Value = -788.5
$ws.Range("H5").Value = 1164
$ws.Range("I5").Value = 960.6667
$ws.Range("J5").Value = 1265.6666
$ws.Range("K5").Value = 2882.0001
$ws.Range("L5").Value = 3796.9998
$ws.Range("M5").Value = -2770.0001
$ws.Range("N5").Value = -4020.9998
$ws.Range("H12").Value = 162.1579
$ws.Range("I12").Value = 72
$ws.Range("J12").Value = 262.33334
$ws.Range("K12").Value = 216
$ws.Range("L12").Value = 787.0000200000001
$ws.Range("M12").Value = -43
$ws.Range("N12").Value = -1133.00002
$ws.Range("H19").Value = 2401
$ws.Range("J19").Value = 2401
$ws.Range("L19").Value = 7203
$ws.Range("N19").Value = -7551
$ws.Range("H35").Value = 396.25
$ws.Range("I35").Value = 592.5
$ws.Range("K35").Value = 1777.5
$ws.Range("M35").Value = -1489.5
$ws.Range("H50").Value = 1866.6666
$ws.Range("I50").Value = 2300
$ws.Range("J50").Value = 1650
$ws.Range("K50").Value = 6900
$ws.Range("L50").Value = 4950
$ws.Range("M50").Value = -6419
$ws.Range("N50").Value = -5912
$ws.Range("H53").Value = 1866.6666
$ws.Range("I53").Value = 2300
$ws.Range("J53").Value = 1650
$ws.Range("K53").Value = 6900
$ws.Range("L53").Value = 4950
$ws.Range("M53").Value = -6419
$ws.Range("N53").Value = -5912
$ws.Range("H87").Value = 3640
$ws.Range("I87").Value = 3640
$ws.Range("K87").Value = 10920
$ws.Range("M87").Value = -9672
$ws.Range("H90").Value = 3640
$ws.Range("I90").Value = 3640
$ws.Range("K90").Value = 32760
$ws.Range("M90").Value = -26520
$ws.Range("H112").Value = 5462.3335
$ws.Range("I112").Value = 6949.5
$ws.Range("J112").Value = 2488
$ws.Range("K112").Value = 20848.5
$ws.Range("L112").Value = 7464
$ws.Range("M112").Value = -19740.5
$ws.Range("N112").Value = -9680
$ws.Range("H116").Value = 15114
$ws.Range("I116").Value = 17304.834
$ws.Range("J116").Value = 1969
$ws.Range("K116").Value = 51914.50199999999
$ws.Range("L116").Value = 5907
$ws.Range("M116").Value = -48472.50199999999
$ws.Range("N116").Value = -12791
$ws.Range("H127").Value = 2500
$ws.Range("J127").Value = 2500
$ws.Range("L127").Value = 7500
$ws.Range("N127").Value = -17420
$ws.Range("H131").Value = 13159722
$ws.Range("J131").Value = 1911.9714
$ws.Range("L131").Value = 5735.914199999999
$ws.Range("N131").Value = -15815.9142
$ws.Range("H134").Value = 5364.7
$ws.Range("I134").Value = 5319.6
$ws.Range("K134").Value = 15958.8
$ws.Range("M134").Value = -10888.8
$ws.Range("H135").Value = 1164
$ws.Range("I135").Value = 960.6667
$ws.Range("J135").Value = 1265.6666
$ws.Range("K135").Value = 8646.0003
$ws.Range("L135").Value = 11390.9994
$ws.Range("M135").Value = -6111.0003
$ws.Range("N135").Value = -16460.9994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 756.1429000000001
$ws.Range("I31").Value = 756.1429000000001
$ws.Range("K31").Value = 756.1429000000001
$ws.Range("M31").Value = -464.1429000000001
$ws.Range("H37").Value = 756.1429000000001
$ws.Range("I37").Value = 756.1429000000001
$ws.Range("K37").Value = 756.1429000000001
$ws.Range("M37").Value = -479.1429000000001
$ws.Range("H70").Value = 4558.245
$ws.Range("I70").Value = 4340.7437
$ws.Range("J70").Value = 5164.143
$ws.Range("K70").Value = 4340.7437
$ws.Range("L70").Value = 5164.143
$ws.Range("M70").Value = -4070.7437
$ws.Range("N70").Value = -5704.143
$ws.Range("H73").Value = 4558.245
$ws.Range("I73").Value = 4340.7437
$ws.Range("J73").Value = 5164.143
$ws.Range("K73").Value = 4340.7437
$ws.Range("L73").Value = 5164.143
$ws.Range("M73").Value = -3404.7437
$ws.Range("N73").Value = -7036.143
$ws.Range("H80").Value = 4520.4194
$ws.Range("I80").Value = 4123.3335
$ws.Range("J80").Value = 4771.2104
$ws.Range("K80").Value = 4123.3335
$ws.Range("L80").Value = 4771.2104
$ws.Range("M80").Value = -3125.3335
$ws.Range("N80").Value = -6767.2104
$ws.Range("H83").Value = 4520.4194
$ws.Range("I83").Value = 4123.3335
$ws.Range("J83").Value = 4771.2104
$ws.Range("K83").Value = 20616.6675
$ws.Range("L83").Value = 23856.052
$ws.Range("M83").Value = -15624.6675
$ws.Range("N83").Value = -33840.052
$ws.Range("H102").Value = 40793.5
$ws.Range("I102").Value = 2428.2273
$ws.Range("K102").Value = 2428.2273
$ws.Range("M102").Value = -806.2273
$ws.Range("H113").Value = 9026.5
$ws.Range("I113").Value = 2880.5
$ws.Range("K113").Value = 2880.5
$ws.Range("M113").Value = -710.5
$ws.Range("H122").Value = 6589.8687
$ws.Range("I122").Value = 5081.1763
$ws.Range("J122").Value = 7811.1904
$ws.Range("K122").Value = 15243.5289
$ws.Range("L122").Value = 23433.5712
$ws.Range("M122").Value = -12793.5289
$ws.Range("N122").Value = -28333.5712
$ws.Range("H126").Value = 9870
$ws.Range("I126").Value = 11515.714
$ws.Range("J126").Value = 4110
$ws.Range("K126").Value = 34547.142
$ws.Range("L126").Value = 12330
$ws.Range("M126").Value = -32077.142
$ws.Range("N126").Value = -17270
$ws.Range("H132").Value = 34051.938
$ws.Range("I132").Value = 38268.62
$ws.Range("K132").Value = 114805.86
$ws.Range("M132").Value = -112275.86
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7584.25
$ws.Range("I7").Value = 7651.1
$ws.Range("J7").Value = 7250
$ws.Range("K7").Value = 7651.1
$ws.Range("L7").Value = 7250
$ws.Range("M7").Value = -7539.1
$ws.Range("N7").Value = -7474
$ws.Range("H16").Value = 2127.6775
$ws.Range("I16").Value = 2176.2727
$ws.Range("J16").Value = 2008.8889
$ws.Range("K16").Value = 2176.2727
$ws.Range("L16").Value = 2008.8889
$ws.Range("M16").Value = -2006.2727
$ws.Range("N16").Value = -2348.8889
$ws.Range("H22").Value = 1327.8462
$ws.Range("I22").Value = 1153
$ws.Range("K22").Value = 1153
$ws.Range("M22").Value = -858
$ws.Range("H27").Value = 1327.8462
$ws.Range("I27").Value = 1153
$ws.Range("K27").Value = 1153
$ws.Range("M27").Value = -1046
$ws.Range("H46").Value = 2034.5294
$ws.Range("I46").Value = 2995.6
$ws.Range("J46").Value = 1634.0834
$ws.Range("K46").Value = 2995.6
$ws.Range("L46").Value = 1634.0834
$ws.Range("M46").Value = -2807.6
$ws.Range("N46").Value = -2010.0834
$ws.Range("H55").Value = 7348.696
$ws.Range("I55").Value = 899.4211
$ws.Range("J55").Value = 37982.75
$ws.Range("K55").Value = 899.4211
$ws.Range("L55").Value = 37982.75
$ws.Range("M55").Value = -726.4211
$ws.Range("N55").Value = -38328.75
$ws.Range("H61").Value = 3232
$ws.Range("I61").Value = 3112.75
$ws.Range("K61").Value = 3112.75
$ws.Range("M61").Value = -2910.75
$ws.Range("H68").Value = 9694.5
$ws.Range("I68").Value = 10147.117
$ws.Range("K68").Value = 10147.117
$ws.Range("M68").Value = -9398.117
$ws.Range("H71").Value = 9694.5
$ws.Range("I71").Value = 10147.117
$ws.Range("K71").Value = 50735.585
$ws.Range("M71").Value = -46991.585
$ws.Range("H93").Value = 50002336
$ws.Range("I93").Value = 2040.4667
$ws.Range("J93").Value = 200003230
$ws.Range("K93").Value = 2040.4667
$ws.Range("L93").Value = 200003230
$ws.Range("M93").Value = -792.4666999999999
$ws.Range("N93").Value = -200005726
$ws.Range("H113").Value = 3232
$ws.Range("I113").Value = 3112.75
$ws.Range("K113").Value = 3112.75
$ws.Range("M113").Value = -942.75
$ws.Range("H122").Value = 3904.375
$ws.Range("I122").Value = 3933.5715
$ws.Range("J122").Value = 3700
$ws.Range("K122").Value = 11800.7145
$ws.Range("L122").Value = 11100
$ws.Range("M122").Value = -9350.7145
$ws.Range("N122").Value = -16000
$ws.Range("H126").Value = 7584.25
$ws.Range("I126").Value = 7651.1
$ws.Range("J126").Value = 7250
$ws.Range("K126").Value = 22953.3
$ws.Range("L126").Value = 21750
$ws.Range("M126").Value = -20483.3
$ws.Range("N126").Value = -26690
$ws.Range("H132").Value = 5594.3887
$ws.Range("I132").Value = 5570.647
$ws.Range("K132").Value = 16711.941
$ws.Range("M132").Value = -14181.941
$ws.Range("H136").Value = 6888.8
$ws.Range("I136").Value = 1111
$ws.Range("K136").Value = 3333
$ws.Range("M136").Value = -783
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 435.9
$ws.Range("I107").Value = 435.9
$ws.Range("K107").Value = 1307.7
$ws.Range("M107").Value = 612.3000000000002
$ws.Range("H113").Value = 1045.7693
$ws.Range("I113").Value = 1060.3
$ws.Range("K113").Value = 3180.9
$ws.Range("M113").Value = -1010.9
$ws.Range("H126").Value = 3811.8667
$ws.Range("I126").Value = 3598.4167
$ws.Range("K126").Value = 10795.2501
$ws.Range("M126").Value = -8325.250100000001
$ws.Range("H132").Value = 1074.2
$ws.Range("I132").Value = 979.06665
$ws.Range("J132").Value = 1359.6
$ws.Range("K132").Value = 2937.19995
$ws.Range("L132").Value = 4078.8
$ws.Range("M132").Value = -407.1999500000002
$ws.Range("N132").Value = -9138.799999999999
$ws.Range("H136").Value = 772930.4
$ws.Range("I136").Value = 1002309.6
$ws.Range("K136").Value = 3006928.8
$ws.Range("M136").Value = -3004378.8
